$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.347
$ws.Range("C4").Value = -11.699
$ws.Range("B6").Value = 5.973000000000001
$ws.Range("B7").Value = 5.374000000000001
$ws.Range("B8").Value = 6.470000000000001
$ws.Range("C8").Value = -12.593
$ws.Range("C9").Value = -11.149
$ws.Range("C12").Value = -10.956
$ws.Range("B16").Value = 4.916
$ws.Range("C17").Value = -13.127
$ws.Range("C18").Value = -11.875
$ws.Range("C19").Value = -11.949
$ws.Range("B20").Value = 8.914999999999999
$ws.Range("C20").Value = -12.068
$ws.Range("B21").Value = 8.580000000000002
$ws.Range("C26").Value = -12.562
$ws.Range("B28").Value = 5.048
$ws.Range("B29").Value = 5.326
$ws.Range("B30").Value = 6.327999999999999
$ws.Range("C31").Value = -12.867
$ws.Range("B32").Value = 6.943
$ws.Range("C39").Value = -12.059
$ws.Range("B40").Value = 9.428999999999998
$ws.Range("C40").Value = -12.206
$ws.Range("C41").Value = -12.02
$ws.Range("C42").Value = -12.364
$ws.Range("C43").Value = -12.395
$ws.Range("B46").Value = 5.298999999999999
$ws.Range("C47").Value = -12.93
$ws.Range("C48").Value = -11.844
$ws.Range("B51").Value = 5.399
$ws.Range("B52").Value = 5.628
$ws.Range("C54").Value = -12.825
$ws.Range("B57").Value = 4.911999999999999
$ws.Range("B59").Value = 5.004
$ws.Range("B62").Value = 5.144
$ws.Range("C62").Value = -13.742
$ws.Range("C63").Value = -10.956
$ws.Range("C64").Value = -11.073
$ws.Range("B66").Value = 4.814
$ws.Range("B73").Value = 6.102000000000001
$ws.Range("B74").Value = 9.186999999999999
$ws.Range("C76").Value = -12.288
$ws.Range("B77").Value = 6.354000000000001
$ws.Range("C81").Value = -13.349
$ws.Range("C84").Value = -13.318
$ws.Range("C89").Value = -13.299
$ws.Range("B92").Value = 4.891
$ws.Range("C94").Value = -11.769
$ws.Range("B100").Value = 6.121
